$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 2452
$ws.Range('K3').Value = 8183
$ws.Range('L3').Value = 2469
$ws.Range('L4').Value = 676
$ws.Range('L5').Value = 147
$ws.Range('L6').Value = 2237
$ws.Range('K7').Value = 27557
$ws.Range('L7').Value = 7981

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L2').Value = 67
$ws.Range('L6').Value = 61
$ws.Range('L7').Value = 262
$ws.Range('L8').Value = 508
$ws.Range('L10').Value = 51
$ws.Range('L11').Value = 140
$ws.Range('L19').Value = 225
$ws.Range('L20').Value = 205
$ws.Range('L22').Value = 26
$ws.Range('L23').Value = 81
$ws.Range('L29').Value = 409
$ws.Range('L31').Value = 76
$ws.Range('L33').Value = 362
$ws.Range('L34').Value = 51
$ws.Range('L37').Value = 292
$ws.Range('L41').Value = 36
$ws.Range('L42').Value = 255
$ws.Range('L43').Value = 65
$ws.Range('L47').Value = 62
$ws.Range('L48').Value = 110
$ws.Range('L50').Value = 43
$ws.Range('L51').Value = 92
$ws.Range('L53').Value = 98
$ws.Range('L55').Value = 72
$ws.Range('K63').Value = 157
$ws.Range('L63').Value = 22
$ws.Range('L65').Value = 147
$ws.Range('L67').Value = 296
$ws.Range('L68').Value = 25
$ws.Range('L73').Value = 62
$ws.Range('L75').Value = 33
$ws.Range('L76').Value = 92
$ws.Range('L77').Value = 49
$ws.Range('L78').Value = 103
$ws.Range('L83').Value = 193
$ws.Range('L84').Value = 80
$ws.Range('L85').Value = 421
$ws.Range('L90').Value = 79
$ws.Range('L94').Value = 95
$ws.Range('L99').Value = 129
$ws.Range('K101').Value = 27557
$ws.Range('L101').Value = 7981

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L2').Value = 75
$ws.Range('L3').Value = 84
$ws.Range('L6').Value = 74
$ws.Range('L7').Value = 262

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('L2').Value = 49
$ws.Range('L3').Value = 44
$ws.Range('L7').Value = 140

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L2').Value = 125
$ws.Range('L3').Value = 170
$ws.Range('L4').Value = 34
$ws.Range('L7').Value = 421

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('L2').Value = 50
$ws.Range('L4').Value = 14

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('L3').Value = 21
$ws.Range('L7').Value = 98

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L2').Value = 146
$ws.Range('L3').Value = 169
$ws.Range('L6').Value = 135
$ws.Range('L7').Value = 508

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('L2').Value = 60
$ws.Range('L7').Value = 193

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L2').Value = 99
$ws.Range('L7').Value = 362

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L2').Value = 87
$ws.Range('L3').Value = 86
$ws.Range('L4').Value = 19
$ws.Range('L7').Value = 292

$ws = $wb.Worksheets.Item('New City')
$ws.Range('L6').Value = 38
$ws.Range('L7').Value = 147

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('L2').Value = 31
$ws.Range('L3').Value = 56
$ws.Range('L6').Value = 30
$ws.Range('L7').Value = 129

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('L3').Value = 18
$ws.Range('L7').Value = 76

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L2').Value = 89
$ws.Range('L4').Value = 25
$ws.Range('L7').Value = 296

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('L2').Value = 29
$ws.Range('L3').Value = 32
$ws.Range('L7').Value = 80

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L2').Value = 134
$ws.Range('L3').Value = 146
$ws.Range('L5').Value = 7
$ws.Range('L6').Value = 106
$ws.Range('L7').Value = 409

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('L4').Value = 28
$ws.Range('L7').Value = 110

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('L3').Value = 70
$ws.Range('L7').Value = 225

$ws = $wb.Worksheets.Item('River North')
$ws.Range('L6').Value = 44
$ws.Range('L7').Value = 92

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('L3').Value = 18
$ws.Range('L7').Value = 61

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('L6').Value = 10
$ws.Range('L7').Value = 36

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('L5').Value = 7
$ws.Range('L7').Value = 255

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('L2').Value = 24
$ws.Range('L7').Value = 51

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('L6').Value = 30
$ws.Range('L7').Value = 103

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('L2').Value = 29
$ws.Range('L4').Value = 4
$ws.Range('L7').Value = 72

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('L2').Value = 26
$ws.Range('L7').Value = 81

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L2').Value = 64
$ws.Range('L3').Value = 64
$ws.Range('L4').Value = 18
$ws.Range('L7').Value = 205

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('L6').Value = 19
$ws.Range('L7').Value = 51

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('L3').Value = 24
$ws.Range('L7').Value = 95

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('L3').Value = 22
$ws.Range('L7').Value = 62

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('L3').Value = 10
$ws.Range('L7').Value = 43

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('L6').Value = 15
$ws.Range('L7').Value = 62

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('L6').Value = 21
$ws.Range('L7').Value = 67

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range('L2').Value = 17
$ws.Range('L3').Value = 12
$ws.Range('L7').Value = 33

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('L6').Value = 19
$ws.Range('L7').Value = 79

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('L2').Value = 27
$ws.Range('L7').Value = 92

$ws = $wb.Worksheets.Item('North Park')
$ws.Range('L4').Value = 3
$ws.Range('L6').Value = 8
$ws.Range('L7').Value = 25

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('L3').Value = 19
$ws.Range('L7').Value = 65

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('L3').Value = 10
$ws.Range('L7').Value = 26

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('L3').Value = 17
$ws.Range('L7').Value = 49
